# Updated cryptos list with latest price/volume snapshot.
# Note: Price cells that look like a plain decimal number (e.g. "4.62")
# are written with a leading apostrophe - PowerShell single-quoted string
# syntax escapes a literal apostrophe as '' - so Excel stores them as text
# (matching the sheet's existing text-formatted Price column) instead of
# silently parsing them into floating-point numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '54.426.02'
$ws.Range("E2").Value = '  +0.56%  '

$ws.Range("D3").Value = '2.289.44'
$ws.Range("E3").Value = '  +1.38%  '

$ws.Range("D4").Value = '''0.998'
$ws.Range("E4").Value = '  -0.35%  '

$ws.Range("D5").Value = '''495.49'
$ws.Range("E5").Value = '  +1.21%  '

$ws.Range("D6").Value = '''127.39'
$ws.Range("E6").Value = '  +0.52%  '

$ws.Range("D7").Value = '''0.997'
$ws.Range("E7").Value = '  -0.29%  '

$ws.Range("E8").Value = '  +1.55%  '

$ws.Range("D9").Value = '2.287.76'
$ws.Range("E9").Value = '  +1.02%  '

$ws.Range("D10").Value = '''0.0945'
$ws.Range("E10").Value = '  +2.47%  '

$ws.Range("E11").Value = '  +2.23%  '

$ws.Range("E12").Value = '  +2.85%  '

$ws.Range("D13").Value = '''4.62'
$ws.Range("E13").Value = '  -2.24%  '

$ws.Range("D14").Value = '2.693.12'
$ws.Range("E14").Value = '  +1.19%  '

$ws.Range("D15").Value = '''21.75'
$ws.Range("E15").Value = '  +2.37%  '

$ws.Range("D16").Value = '54.161.60'
$ws.Range("E16").Value = '  +0.16%  '

$ws.Range("E17").Value = '  +0.54%  '

$ws.Range("D18").Value = '2.281.41'
$ws.Range("E18").Value = '  +1.42%  '

$ws.Range("D19").Value = '''10.06'
$ws.Range("E19").Value = '  +4.49%  '

$ws.Range("D20").Value = '''4.09'
$ws.Range("E20").Value = '  +2.63%  '

$ws.Range("B21").Value = 'Uniswap'
$ws.Range("C21").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D21").Value = '''6.46'
$ws.Range("E21").Value = '  +5.67%  '

$ws.Range("B22").Value = 'BitcoinCash'
$ws.Range("C22").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D22").Value = '''301.53'
$ws.Range("E22").Value = '  +0.16%  '

$ws.Range("E23").Value = '  +0.02%  '

$ws.Range("E24").Value = '  -1.76%  '

$ws.Range("D25").Value = '''62.31'
$ws.Range("E25").Value = '  -2.45%  '

$ws.Range("E26").Value = '  -0.33%  '

$ws.Range("E27").Value = '  +1.41%  '

$ws.Range("D28").Value = '2.394.55'
$ws.Range("E28").Value = '  +1.26%  '

$ws.Range("E29").Value = '  +4.62%  '

$ws.Range("E30").Value = '  +0.07%  '

$ws.Range("D31").Value = '''168.55'
$ws.Range("E31").Value = '  -0.98%  '

$ws.Range("E32").Value = '  -0.53%  '

$ws.Range("E33").Value = '  -0.19%  '

$ws.Range("D34").Value = '''5.86'
$ws.Range("E34").Value = '  +1.82%  '

$ws.Range("D35").Value = '''0.998'
$ws.Range("E35").Value = '  -0.07%  '

$ws.Range("E36").Value = '  +0.09%  '

$ws.Range("E37").Value = '  +0.37%  '

$ws.Range("D38").Value = '''17.56'
$ws.Range("E38").Value = '  +0.44%  '

$ws.Range("E39").Value = '  +2.12%  '

$ws.Range("D40").Value = '''0.868'
$ws.Range("E40").Value = '  +2.12%  '

$ws.Range("E41").Value = '  +3.23%  '

$ws.Range("D42").Value = '''35.38'
$ws.Range("E42").Value = '  -1.14%  '

$ws.Range("E43").Value = '  +1.86%  '

$ws.Range("E44").Value = '  +1.50%  '

$ws.Range("D45").Value = '''3.34'

$ws.Range("B46").Value = 'Aave'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D46").Value = '''128.55'
$ws.Range("E46").Value = '  +5.11%  '

$ws.Range("B47").Value = 'RenderToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D47").Value = '''4.79'
$ws.Range("E47").Value = '  +2.10%  '

$ws.Range("D48").Value = '''0.0889'
$ws.Range("E48").Value = '  +0.93%  '

$ws.Range("E49").Value = '  +0.47%  '

$ws.Range("D50").Value = '''238.00'
$ws.Range("E50").Value = '  +0.08%  '

$ws.Range("E51").Value = '  +2.28%  '
